$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 487
$ws.Range("A487").Value = 46
$ws.Range("B487").Value = "Simulator_PE"
$ws.Range("C487").Value = "Control"
$ws.Range("D487").Value = "diff"
$ws.Range("E487").Value = "YP_Total_points"
$ws.Range("F487").Value = "continuous"
$ws.Range("G487").Value = 5.24831
$ws.Range("H487").Value = 0.78194
$ws.Range("I487").Value = 0.99879
$ws.Range("J487").Value = 5.09173
$ws.Range("K487").Value = 0.83806
$ws.Range("L487").Value = 4.96192
$ws.Range("M487").Value = 0.87127
$ws.Range("N487").Value = 0.8705558429412465
$ws.Range("O487").Value = 0.1931927546610602
$ws.Range("P487").Value = 0.6129125145199258
$ws.Range("Q487").Value = 1.080824786354194
$ws.Range("R487").Value = -0.1518554135797881

# Row 488
$ws.Range("A488").Value = 46
$ws.Range("B488").Value = "Simulator_PE"
$ws.Range("C488").Value = "Control"
$ws.Range("D488").Value = "diff"
$ws.Range("E488").Value = "YS_Patient_interaction_points"
$ws.Range("F488").Value = "continuous"
$ws.Range("G488").Value = -0.09515
$ws.Range("H488").Value = 0.1068
$ws.Range("I488").Value = 0.13869
$ws.Range("J488").Value = -0.06962
$ws.Range("K488").Value = 0.1161
$ws.Range("L488").Value = -0.04236
$ws.Range("M488").Value = 0.12088
$ws.Range("N488").Value = 0.8462098298045659
$ws.Range("O488").Value = -0.2288722101142771
$ws.Range("P488").Value = 0.592996607871869
$ws.Range("Q488").Value = 1.084037892428414
$ws.Range("R488").Value = 0.2300148227740146

# Row 489
$ws.Range("A489").Value = 46
$ws.Range("B489").Value = "Simulator_PE"
$ws.Range("C489").Value = "Control"
$ws.Range("D489").Value = "diff"
$ws.Range("E489").Value = "YS_Inspection_points"
$ws.Range("F489").Value = "continuous"
$ws.Range("G489").Value = 1.87892
$ws.Range("H489").Value = 0.19921
$ws.Range("I489").Value = 0.25278
$ws.Range("J489").Value = 1.78624
$ws.Range("K489").Value = 0.22452
$ws.Range("L489").Value = 1.76782
$ws.Range("M489").Value = 0.22832
$ws.Range("N489").Value = 0.7872491562577
$ws.Range("O489").Value = 0.4366700833511567
$ws.Range("P489").Value = 0.6210647090859649
$ws.Range("Q489").Value = 1.034136446566341
$ws.Range("R489").Value = -0.08135037317159913

# Row 490
$ws.Range("A490").Value = 46
$ws.Range("B490").Value = "Simulator_PE"
$ws.Range("C490").Value = "Control"
$ws.Range("D490").Value = "diff"
$ws.Range("E490").Value = "YS_Pulses_points"
$ws.Range("F490").Value = "continuous"
$ws.Range("G490").Value = 1.35759
$ws.Range("H490").Value = 0.29057
$ws.Range("I490").Value = 0.36923
$ws.Range("J490").Value = 1.39476
$ws.Range("K490").Value = 0.31145
$ws.Range("L490").Value = 1.2989
$ws.Range("M490").Value = 0.32481
$ws.Range("N490").Value = 0.8704120131005727
$ws.Range("O490").Value = -0.1234100653056829
$ws.Range("P490").Value = 0.6193092778528101
$ws.Range("Q490").Value = 1.087632340055104
$ws.Range("R490").Value = -0.3012569530522705

# Row 491
$ws.Range("A491").Value = 46
$ws.Range("B491").Value = "Simulator_PE"
$ws.Range("C491").Value = "Control"
$ws.Range("D491").Value = "diff"
$ws.Range("E491").Value = "YS_Auscultation_points"
$ws.Range("F491").Value = "continuous"
$ws.Range("G491").Value = 1.3297
$ws.Range("H491").Value = 0.46615
$ws.Range("I491").Value = 0.59279
$ws.Range("J491").Value = 1.26126
$ws.Range("K491").Value = 0.4791
$ws.Range("L491").Value = 1.20577
$ws.Range("M491").Value = 0.48839
$ws.Range("N491").Value = 0.9466709174010041
$ws.Range("O491").Value = 0.1447946640101536
$ws.Range("P491").Value = 0.6183717482416357
$ws.Range("Q491").Value = 1.039157040215031
$ws.Range("R491").Value = -0.1147039081406277

# Row 492
$ws.Range("A492").Value = 46
$ws.Range("B492").Value = "Simulator_PE"
$ws.Range("C492").Value = "Control"
$ws.Range("D492").Value = "diff"
$ws.Range("E492").Value = "YS_Cardiac_findings_points"
$ws.Range("F492").Value = "continuous"
$ws.Range("G492").Value = 0.5289700000000001
$ws.Range("H492").Value = 0.18294
$ws.Range("I492").Value = 0.23943
$ws.Range("J492").Value = 0.46519
$ws.Range("K492").Value = 0.19692
$ws.Range("L492").Value = 0.46841
$ws.Range("M492").Value = 0.19555
$ws.Range("N492").Value = 0.8630534546457863
$ws.Range("O492").Value = 0.3355807402565573
$ws.Range("P492").Value = 0.5837947948078227
$ws.Range("Q492").Value = 0.986134121885734
$ws.Range("R492").Value = 0.01640879752444094

# Row 493
$ws.Range("A493").Value = 47
$ws.Range("B493").Value = "HFNC"
$ws.Range("C493").Value = "NIV"
$ws.Range("D493").Value = "diff"
$ws.Range("E493").Value = "YP_intubation"
$ws.Range("F493").Value = "binary"
$ws.Range("G493").Value = 0.0772
$ws.Range("H493").Value = 0.04141
$ws.Range("I493").Value = 0.05336
$ws.Range("J493").Value = -0.06349
$ws.Range("K493").Value = 0.05547
$ws.Range("L493").Value = 0.06894
$ws.Range("M493").Value = 0.04103
$ws.Range("N493").Value = 0.5573065267943678
$ws.Range("O493").Value = 2.874306054592232
$ws.Range("P493").Value = 0.6022527880549982
$ws.Range("Q493").Value = 0.5471251795910645
$ws.Range("R493").Value = 2.714441483635908

# Row 494
$ws.Range("A494").Value = 47
$ws.Range("B494").Value = "HFNC"
$ws.Range("C494").Value = "NIV"
$ws.Range("D494").Value = "diff"
$ws.Range("E494").Value = "YS_nasal_pressure_injury"
$ws.Range("F494").Value = "binary"
$ws.Range("G494").Value = -0.10489
$ws.Range("H494").Value = 0.02362
$ws.Range("I494").Value = 0.03109
$ws.Range("J494").Value = -0.0873
$ws.Range("K494").Value = 0.02525
$ws.Range("L494").Value = -0.09209000000000001
$ws.Range("M494").Value = 0.02535
$ws.Range("N494").Value = 0.8750583668267814
$ws.Range("O494").Value = -0.7194689556885784
$ws.Range("P494").Value = 0.5771893943242214
$ws.Range("Q494").Value = 1.007936476815998
$ws.Range("R494").Value = -0.1893276935128095

# Row 495
$ws.Range("A495").Value = 47
$ws.Range("B495").Value = "HFNC"
$ws.Range("C495").Value = "NIV"
$ws.Range("D495").Value = "diff"
$ws.Range("E495").Value = "YS_antibiotic"
$ws.Range("F495").Value = "binary"
$ws.Range("G495").Value = 0.0512
$ws.Range("H495").Value = 0.05595
$ws.Range("I495").Value = 0.07212
$ws.Range("J495").Value = -0.09762999999999999
$ws.Range("K495").Value = 0.0631
$ws.Range("L495").Value = -0.09143999999999999
$ws.Range("M495").Value = 0.06071
$ws.Range("N495").Value = 0.7862152496100824
$ws.Range("O495").Value = 2.495796802325497
$ws.Range("P495").Value = 0.6018506662495396
$ws.Range("Q495").Value = 0.925681847292929
$ws.Range("R495").Value = 0.09997329803747372

